$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.158.00'
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").Value = '1.902.10'
$ws.Range("E3").Value = '  +0.24%  '
$ws.Range("E4").Value = '  -0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '252.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.700'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.95%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.85'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.34%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.356'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.37'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.15%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0759'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.10%  '
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '13.22'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.07%  '
$ws.Range("D14").Value = '2.177.15'
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("E15").Value = '  +3.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.01'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.64%  '
$ws.Range("D17").Value = '1.906.97'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '35.157.99'
$ws.Range("E18").Value = '  -0.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.03%  '
$ws.Range("D20").Value = '0.0₃0843'
$ws.Range("E20").Value = '  +2.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '243.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.09'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.85%  '
$ws.Range("E24").Value = '  -0.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.42'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.49%  '
$ws.Range("E26").Value = '  -1.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.18'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.59'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.61%  '
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.130'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.00%  '
$ws.Range("D31").Value = '4.128.38'
$ws.Range("E31").Value = '  -0.34%  '
$ws.Range("E32").Value = '  +12.11%  '
$ws.Range("E33").Value = '  +4.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0595'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.64%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +7.48%  '
$ws.Range("E36").Value = '  +3.34%  '
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.845'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.46'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.15%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.18'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0216'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0670'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.79%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.11'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.27%  '
$ws.Range("D46").Value = '1.306.40'
$ws.Range("E46").Value = '  -3.29%  '
$ws.Range("E47").Value = '  -0.18%  '
$ws.Range("E48").Value = '  -1.38%  '
$ws.Range("E49").Value = '  +1.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '12.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.64%  '
$ws.Range("E51").Value = '  +6.79%  '
